# Rename the worksheet to match the player name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Jos Buttler"

# New header row, with an extra leading "matchNo" column.
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

# Full set of innings rows (matchNo, teamName, batterName, states, runs, balls, fours, sixes, sr, opponentTeamName, venue, date, result).
$rows = @(
    @("28th","Rajasthan Royals","Jos Buttler","b Sandeep Sharma","124","64","11","8","193.75","Sunrisers Hyderabad","Delhi","May 02","Royals won by 55 runs"),
    @("24th","Rajasthan Royals","Jos Buttler","st †de Kock b Chahar","41","32","3","3","128.12","Mumbai Indians","Delhi","April 29","Mumbai won by 7 wickets (with 9 balls remaining)"),
    @("18th","Rajasthan Royals","Jos Buttler","lbw b Varun","5","7","1","0","71.42","Kolkata Knight Riders","Wankhede","April 24","Royals won by 6 wickets (with 7 balls remaining)"),
    @("16th","Rajasthan Royals","Jos Buttler","b Mohammed Siraj","8","8","2","0","100.00","Royal Challengers Bangalore","Wankhede","April 22","RCB won by 10 wickets (with 21 balls remaining)"),
    @("4th","Rajasthan Royals","Jos Buttler","b Richardson","25","13","5","0","192.30","Punjab Kings","Wankhede","April 12","Punjab Kings won by 4 runs"),
    @("7th","Rajasthan Royals","Jos Buttler","c †Pant b Woakes","2","7","0","0","28.57","Delhi Capitals","Wankhede","April 15","Royals won by 3 wickets (with 2 balls remaining)"),
    @("12th","Rajasthan Royals","Jos Buttler","b Jadeja","49","35","5","2","140.00","Chennai Super Kings","Wankhede","April 19","Super Kings won by 45 runs")
)

# Write header row 1.
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Write data rows starting at row 2.
for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
